$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2196.599662171774
$ws.Range("C2").Value = 2205.423038560061
$ws.Range("D2").Value = 2213.267102976351
$ws.Range("E2").Value = 2211.062976753502
$ws.Range("B3").Value = 2201.606885761249
$ws.Range("C3").Value = 2208.593229574046
$ws.Range("D3").Value = 2214.879550218847
$ws.Range("E3").Value = 2212.289353747457
$ws.Range("B4").Value = 2178.408688171938
$ws.Range("C4").Value = 2192.085558436731
$ws.Range("D4").Value = 2204.643962329801
$ws.Range("E4").Value = 2204.188048573605
$ws.Range("B5").Value = 2198.951043606867
$ws.Range("C5").Value = 2206.043190708929
$ws.Range("D5").Value = 2213.264139928402
$ws.Range("E5").Value = 2212.0045224477
$ws.Range("B6").Value = 2206.670585597019
$ws.Range("C6").Value = 2213.253418125554
$ws.Range("D6").Value = 2219.368179815974
$ws.Range("E6").Value = 2216.909327032102
$ws.Range("B7").Value = 2191.721997665262
$ws.Range("C7").Value = 2200.677181678597
$ws.Range("D7").Value = 2209.792567099627
$ws.Range("E7").Value = 2209.28319458538
$ws.Range("B8").Value = 2188.904419689039
$ws.Range("C8").Value = 2196.697326503753
$ws.Range("D8").Value = 2204.063040089982
$ws.Range("E8").Value = 2202.028990745872
$ws.Range("B9").Value = 2194.535756323377
$ws.Range("C9").Value = 2200.595623566658
$ws.Range("D9").Value = 2207.328665015756
$ws.Range("E9").Value = 2205.775692409016
$ws.Range("B10").Value = 2045.207259279136
$ws.Range("C10").Value = 2097.389136720578
$ws.Range("D10").Value = 2152.809230531152
$ws.Range("E10").Value = 2172.380610578276
$ws.Range("B11").Value = 2033.936097627198
$ws.Range("C11").Value = 2104.203142466426
$ws.Range("D11").Value = 2155.703348248459
$ws.Range("E11").Value = 2173.445557392237
$ws.Range("B12").Value = 1876.610898948367
$ws.Range("C12").Value = 1991.94078381052
$ws.Range("D12").Value = 2098.051720771348
$ws.Range("E12").Value = 2134.639480504867
$ws.Range("B13").Value = 2034.261254671876
$ws.Range("C13").Value = 2087.825279412683
$ws.Range("D13").Value = 2145.920377647311
$ws.Range("E13").Value = 2167.399040879289
